# This workbook ("Rspo3-Lrp6.xlsx") recomputed its NATMI TPM-based ligand-receptor
# statistics. The "Sending cluster" / "Target cluster" labels (columns A and D) and
# the "Ligand-expressing cells" / "Receptor-expressing cells" counts (columns E, K, L)
# are unchanged; only the TPM-derived expression/specificity/weight numbers in columns
# F-J, M-T were refreshed with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => hashtable of column letter => new numeric value
$updates = @{
    2  = @{ G = 0.06762866666666667; H = 0.202886;   I = 0.0134153952845566; J = 0.0134153952845566;
            M = 9.841031333333333;   N = 29.523094;  O = 0.1083017349730097; P = 0.1125970533891552;
            Q = 0.6655358276982223;  R = 5.989822449284; S = 0.001452910584666213; T = 0.00151053397909184 }
    3  = @{ G = 0.06762866666666667; H = 0.202886;   I = 0.0134153952845566; J = 0.0134153952845566;
            O = 0.3504595127507141;  P = 0.3643589687437936;
            Q = 2.153643816984;      R = 19.382794352856; S = 0.004701552894783933; T = 0.004888019591171394 }
    4  = @{ G = 0.06762866666666667; H = 0.202886;   I = 0.0134153952845566; J = 0.0134153952845566;
            M = 16.16670066666667;   N = 48.500102;  O = 0.1779164877830196; P = 0.1849727733236046;
            Q = 1.093332410485778;   R = 9.839991694371999; S = 0.002386820011249192; T = 0.002481482871016842 }
    5  = @{ G = 0.06762866666666667; H = 0.202886;   I = 0.0134153952845566; J = 0.0134153952845566;
            M = 10.399077;           N = 20.798154;  O = 0.1144430947397913; P = 0.07932132236322763;
            Q = 0.7032757120740001;  R = 4.219654272444; S = 0.00153529935352226; T = 0.001064126893996438 }
    6  = @{ G = 0.06762866666666667; H = 0.202886;   I = 0.0134153952845566; J = 0.0134153952845566;
            M = 22.614852;           N = 67.84455600000001; O = 0.2488791697534654; P = 0.258749882180219;
            Q = 1.529412287624;      R = 13.764710588616; S = 0.003338812440335001; T = 0.003471231949280085 }
    7  = @{ I = 0.9827953701592058;  J = 0.9827953701592059;
            M = 9.841031333333333;   N = 29.523094;  O = 0.1083017349730097; P = 0.1125970533891552;
            Q = 48.75633675064733;   R = 438.807030755826; S = 0.1064384437116833; T = 0.1106598627644306 }
    8  = @{ I = 0.9827953701592058;  J = 0.9827953701592059;
            O = 0.3504595127507141;  P = 0.3643589687437936;
            S = 0.344429986559653;   T = 0.3580903075573832 }
    9  = @{ I = 0.9827953701592058;  J = 0.9827953701592059;
            M = 16.16670066666667;   N = 48.500102;  O = 0.1779164877830196; P = 0.1849727733236046;
            Q = 80.09618861602867;   R = 720.865697544258; S = 0.1748555004681385; T = 0.1817903852279469 }
    10 = @{ I = 0.9827953701592058;  J = 0.9827953701592059;
            M = 10.399077;           N = 20.798154;  O = 0.1144430947397913; P = 0.07932132236322763;
            Q = 51.52111429526099;   R = 309.126685771566; S = 0.1124741436569582; T = 0.07795662837348599 }
    11 = @{ I = 0.9827953701592058;  J = 0.9827953701592059;
            M = 22.614852;           N = 67.84455600000001; O = 0.2488791697534654; P = 0.258749882180219;
            Q = 112.042864444836;    R = 1008.385780003524; S = 0.2445972957627729; T = 0.2542981862359592 }
    12 = @{ F = 0.3333333333333333;  G = 0.019102;   H = 0.057306; I = 0.003789234556237495; J = 0.003789234556237496;
            M = 9.841031333333333;   N = 29.523094;  O = 0.1083017349730097; P = 0.1125970533891552;
            Q = 0.1879833805293333;  R = 1.691850424764; S = 0.0004103806766602032; T = 0.000426656645632705 }
    13 = @{ F = 0.3333333333333333;  G = 0.019102;   H = 0.057306; I = 0.003789234556237495; J = 0.003789234556237496;
            O = 0.3504595127507141;  P = 0.3643589687437936;
            Q = 0.6083057114639999;  R = 5.474751403176; S = 0.001327973296277161; T = 0.001380641595239041 }
    14 = @{ F = 0.3333333333333333;  G = 0.019102;   H = 0.057306; I = 0.003789234556237495; J = 0.003789234556237496;
            M = 16.16670066666667;   N = 48.500102;  O = 0.1779164877830196; P = 0.1849727733236046;
            Q = 0.3088163161346667;  R = 2.779346845212; S = 0.0006741673036318238; T = 0.0007009052246408879 }
    15 = @{ F = 0.3333333333333333;  G = 0.019102;   H = 0.057306; I = 0.003789234556237495; J = 0.003789234556237496;
            M = 10.399077;           N = 20.798154;  O = 0.1144430947397913; P = 0.07932132236322763;
            Q = 0.198643168854;      R = 1.191859013124; S = 0.0004336517293107785; T = 0.0003005670957451962 }
    16 = @{ F = 0.3333333333333333;  G = 0.019102;   H = 0.057306; I = 0.003789234556237495; J = 0.003789234556237496;
            M = 22.614852;           N = 67.84455600000001; O = 0.2488791697534654; P = 0.258749882180219;
            Q = 0.431988902904;      R = 3.887900126136001; S = 0.0009430615503575288; T = 0.0009804639949796664 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
